$d = $word.ActiveDocument

# Find the paragraph index of "Ver no Jupiter Salvar em pdf Salvar em docx"
# and the paragraph index of the copyright/footer line that follows it.
$count = $d.Paragraphs.Count
$idx1 = -1
$idx2 = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($idx1 -eq -1 -and $t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $idx1 = $i
    }
    if ($idx2 -eq -1 -and $t -like "*Contact: luizeleno@usp.br*") {
        $idx2 = $i
    }
}

if ($idx1 -ne -1 -and $idx2 -ne -1) {
    # Also remove the blank paragraph immediately preceding "Ver no Jupiter...",
    # which separated it from the "Requisitos" section above.
    $startPara = $d.Paragraphs.Item($idx1 - 1)
    $endPara = $d.Paragraphs.Item($idx2)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
